$d = $word.ActiveDocument

# 1) Update the text of the "O Patrimônio Líquido, como indicador relevante..." paragraph.
$old1 = "O Patrimônio Líquido, como indicador relevante, desempenha um papel crucial na determinação da saúde financeira de uma organização. Quando o Patrimônio Líquido é positivo, isso pode indicar uma base sólida para o crescimento e estabilidade. Por outro lado, um Patrimônio Líquido negativo pode ser um sinal de alerta."
$new1 = "O Patrimônio Líquido, como indicador relevante, desempenha um papel crucial na determinação da saúde financeira de uma organização. É importante destacar que esse valor representa a participação dos sócios ou acionistas na empresa, excluindo o capital investido. Quando o Patrimônio Líquido é positivo, isso pode indicar uma base sólida para o crescimento e estabilidade dos lucros dos sócios. Por outro lado, um Patrimônio Líquido negativo pode ser um sinal de alerta para os sócios, indicando que os passivos superam os ativos e, portanto, pode ser necessário tomar medidas para reequilibrar as finanças."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2) Update the text of the "Este sistema proporciona..." paragraph.
$old2 = "Este sistema proporciona informações importantes para tomadas de decisões estratégicas. Com sua capacidade de gerar relatórios, os usuários podem planejar e implementar medidas para alcançar seus objetivos financeiros."
$new2 = "Este sistema proporciona informações importantes para tomadas de decisões estratégicas. Com sua capacidade de gerar relatórios, os usuários podem planejar e implementar medidas para alcançar seus objetivos financeiros, permitindo que os sócios compreendam melhor a situação financeira da empresa."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) Remove the now-redundant blank paragraph immediately preceding the "O Patrimônio Líquido..." paragraph,
#    and the blank paragraph immediately following the "Este sistema..." paragraph.
#    Locate them by walking paragraphs and matching neighbors (robust to index shifts from steps above,
#    since Find/Replace above did not add/remove paragraphs).
$target1 = $new1.Substring(0, 40)
$target2 = $new2.Substring(0, 40)

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith($target2)) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim().Length -eq 0) {
            $next.Range.Delete()
        }
    }
    if ($t.StartsWith($target1)) {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text.Trim().Length -eq 0) {
            $prev.Range.Delete()
        }
    }
}
